$d = $word.ActiveDocument

$pairs = @(
    @("616÷4=", "827÷9="),
    @("656÷3=", "601÷2="),
    @("433÷8=", "260÷3="),
    @("398÷7=", "425÷3="),
    @("572÷3=", "612÷2="),
    @("486÷2=", "344÷4="),
    @("217÷5=", "819÷6="),
    @("428÷6=", "980÷2="),
    @("584÷2=", "315÷2="),
    @("131÷9=", "458÷5="),
    @("369÷8=", "497÷8="),
    @("352÷3=", "743÷5="),
    @("450÷6=", "789÷9="),
    @("968÷9=", "270÷4="),
    @("916÷9=", "176÷8="),
    @("442÷7=", "953÷2="),
    @("165÷3=", "364÷6="),
    @("531÷2=", "796÷7="),
    @("745÷2=", "347÷7="),
    @("521÷8=", "489÷4="),
    @("907÷4=", "815÷6="),
    @("718÷2=", "648÷3="),
    @("638÷5=", "914÷2="),
    @("518÷3=", "186÷6="),
    @("935÷9=", "586÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
